# Scheduled price-data refresh: update computed Leve profit columns (H:N)
# across all crafting-job sheets with freshly-fetched Market Board averages.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 845.8033
$ws.Range("I15").Value = 845.8033
$ws.Range("K15").Value = 2537.4099
$ws.Range("M15").Value = -2368.4099
# Row 33
$ws.Range("H33").Value = 774.3333
$ws.Range("I33").Value = 330
$ws.Range("K33").Value = 330
$ws.Range("M33").Value = -101
# Row 86
$ws.Range("H86").Value = 2184.625
$ws.Range("J86").Value = 2490.75
$ws.Range("L86").Value = 2490.75
$ws.Range("N86").Value = -4736.75
# Row 89
$ws.Range("H89").Value = 2184.625
$ws.Range("J89").Value = 2490.75
$ws.Range("L89").Value = 12453.75
$ws.Range("N89").Value = -23685.75
# Row 111
$ws.Range("H111").Value = 10561.969
$ws.Range("J111").Value = 5404.6665
$ws.Range("L111").Value = 16213.9995
$ws.Range("N111").Value = -22347.9995
# Row 125
$ws.Range("H125").Value = 2473.75
$ws.Range("I125").Value = 1258
$ws.Range("K125").Value = 11322
$ws.Range("M125").Value = -8862
# Row 135
$ws.Range("H135").Value = 6011.625
$ws.Range("I135").Value = 2912.3333
$ws.Range("K135").Value = 26210.9997
$ws.Range("M135").Value = -23675.9997

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6561.415
$ws.Range("I61").Value = 5642.8
$ws.Range("J61").Value = 7381.607
$ws.Range("K61").Value = 5642.8
$ws.Range("L61").Value = 7381.607
$ws.Range("M61").Value = -5430.8
$ws.Range("N61").Value = -7805.607
# Row 122
$ws.Range("H122").Value = 4131.457
$ws.Range("I122").Value = 3642.4443
$ws.Range("J122").Value = 5781.875
$ws.Range("K122").Value = 10927.3329
$ws.Range("L122").Value = 17345.625
$ws.Range("M122").Value = -8477.332900000001
$ws.Range("N122").Value = -22245.625
# Row 132
$ws.Range("H132").Value = 1768
$ws.Range("I132").Value = 1707.1666
$ws.Range("K132").Value = 5121.4998
$ws.Range("M132").Value = -2591.4998
# Row 134
$ws.Range("H134").Value = 95827.62
$ws.Range("J134").Value = 95479.914
$ws.Range("L134").Value = 95479.914
$ws.Range("N134").Value = -105619.914
# Row 136
$ws.Range("H136").Value = 6561.415
$ws.Range("I136").Value = 5642.8
$ws.Range("J136").Value = 7381.607
$ws.Range("K136").Value = 16928.4
$ws.Range("L136").Value = 22144.821
$ws.Range("M136").Value = -14378.4
$ws.Range("N136").Value = -27244.821

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2566340.5
$ws.Range("I86").Value = 5558372
$ws.Range("K86").Value = 5558372
$ws.Range("M86").Value = -5557249
# Row 89
$ws.Range("H89").Value = 2566340.5
$ws.Range("I89").Value = 5558372
$ws.Range("K89").Value = 27791860
$ws.Range("M89").Value = -27786244
# Row 105
$ws.Range("H105").Value = 5069.273
$ws.Range("I105").Value = 4491.6
$ws.Range("J105").Value = 6874.5
$ws.Range("K105").Value = 4491.6
$ws.Range("L105").Value = 6874.5
$ws.Range("M105").Value = -2744.6
$ws.Range("N105").Value = -10368.5
# Row 134
$ws.Range("H134").Value = 6184.3438
$ws.Range("I134").Value = 2752.2415
$ws.Range("K134").Value = 8256.7245
$ws.Range("M134").Value = -5721.7245

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 333.0909
$ws.Range("I7").Value = 210.42857
$ws.Range("J7").Value = 547.75
$ws.Range("K7").Value = 210.42857
$ws.Range("L7").Value = 547.75
$ws.Range("M7").Value = -97.42857000000001
$ws.Range("N7").Value = -773.75
# Row 16
$ws.Range("H16").Value = 985.8
$ws.Range("J16").Value = 925
$ws.Range("L16").Value = 925
$ws.Range("N16").Value = -1499
# Row 31
$ws.Range("H31").Value = 2595.5898
$ws.Range("I31").Value = 1973.0278
$ws.Range("J31").Value = 10066.333
$ws.Range("K31").Value = 1973.0278
$ws.Range("L31").Value = 10066.333
$ws.Range("M31").Value = -1678.0278
$ws.Range("N31").Value = -10656.333
# Row 34
$ws.Range("H34").Value = 2595.5898
$ws.Range("I34").Value = 1973.0278
$ws.Range("J34").Value = 10066.333
$ws.Range("K34").Value = 1973.0278
$ws.Range("L34").Value = 10066.333
$ws.Range("M34").Value = -1771.0278
$ws.Range("N34").Value = -10470.333
# Row 58
$ws.Range("H58").Value = 1655.6666
$ws.Range("J58").Value = 2479.2
$ws.Range("L58").Value = 2479.2
$ws.Range("N58").Value = -2885.2
# Row 107
$ws.Range("H107").Value = 1400.6086
$ws.Range("I107").Value = 936.2353000000001
$ws.Range("K107").Value = 936.2353000000001
$ws.Range("M107").Value = 983.7646999999999
# Row 113
$ws.Range("H113").Value = 985.8
$ws.Range("J113").Value = 925
$ws.Range("L113").Value = 925
$ws.Range("N113").Value = -5265
# Row 122
$ws.Range("H122").Value = 3964.36
$ws.Range("I122").Value = 3909.4375
$ws.Range("J122").Value = 4062
$ws.Range("K122").Value = 11728.3125
$ws.Range("L122").Value = 12186
$ws.Range("M122").Value = -9278.3125
$ws.Range("N122").Value = -17086
# Row 136
$ws.Range("H136").Value = 1655.6666
$ws.Range("J136").Value = 2479.2
$ws.Range("L136").Value = 7437.599999999999
$ws.Range("N136").Value = -12537.6
# Row 141
$ws.Range("H141").Value = 331793.94
$ws.Range("J141").Value = 349913.53
$ws.Range("L141").Value = 349913.53
$ws.Range("N141").Value = -360273.53

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 425.25
$ws.Range("I86").Value = 459
$ws.Range("K86").Value = 1377
$ws.Range("M86").Value = -191
# Row 89
$ws.Range("H89").Value = 425.25
$ws.Range("I89").Value = 459
$ws.Range("K89").Value = 4131
$ws.Range("M89").Value = 1797
# Row 114
$ws.Range("H114").Value = 4483
$ws.Range("J114").Value = 4483
$ws.Range("L114").Value = 13449
$ws.Range("N114").Value = -19957
# Row 121
$ws.Range("H121").Value = 4516.9
$ws.Range("J121").Value = 4935.4443
$ws.Range("L121").Value = 14806.3329
$ws.Range("N121").Value = -17426.3329

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
# Row 102
$ws.Range("H102").Value = 17008.273
$ws.Range("I102").Value = 1567.28
$ws.Range("K102").Value = 1567.28
$ws.Range("M102").Value = 54.72000000000003
# Row 113
$ws.Range("H113").Value = 2606.7778
$ws.Range("I113").Value = 1921.2858
$ws.Range("K113").Value = 1921.2858
$ws.Range("M113").Value = 248.7141999999999
# Row 122
$ws.Range("H122").Value = 4415.4165
$ws.Range("I122").Value = 3029.6
$ws.Range("J122").Value = 5405.2856
$ws.Range("K122").Value = 9088.799999999999
$ws.Range("L122").Value = 16215.8568
$ws.Range("M122").Value = -6638.799999999999
$ws.Range("N122").Value = -21115.8568
# Row 132
$ws.Range("H132").Value = 2572.1177
$ws.Range("I132").Value = 2305.3655
$ws.Range("K132").Value = 6916.0965
$ws.Range("M132").Value = -4386.0965

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2702.2856
$ws.Range("I16").Value = 1348
$ws.Range("J16").Value = 3244
$ws.Range("K16").Value = 1348
$ws.Range("L16").Value = 3244
$ws.Range("M16").Value = -1178
$ws.Range("N16").Value = -3584
# Row 40
$ws.Range("H40").Value = 4934.263
$ws.Range("J40").Value = 4949.25
$ws.Range("L40").Value = 4949.25
$ws.Range("N40").Value = -5221.25
# Row 55
$ws.Range("H55").Value = 409.70587
$ws.Range("I55").Value = 333.63635
$ws.Range("K55").Value = 333.63635
$ws.Range("M55").Value = -160.63635
# Row 132
$ws.Range("H132").Value = 3182.4424
$ws.Range("I132").Value = 3190.225
$ws.Range("K132").Value = 9570.674999999999
$ws.Range("M132").Value = -7040.674999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 90000
$ws.Range("J16").Value = 90000
$ws.Range("L16").Value = 90000
$ws.Range("N16").Value = -90584
# Row 96
$ws.Range("H96").Value = 3230.375
$ws.Range("I96").Value = 2260
$ws.Range("K96").Value = 2260
$ws.Range("M96").Value = -887
# Row 113
$ws.Range("H113").Value = 3475356.2
$ws.Range("I113").Value = 3971803.5
$ws.Range("J113").Value = 225.33333
$ws.Range("K113").Value = 11915410.5
$ws.Range("L113").Value = 675.99999
$ws.Range("M113").Value = -11913240.5
$ws.Range("N113").Value = -5015.99999
# Row 122
$ws.Range("H122").Value = 2929.2144
$ws.Range("I122").Value = 3024
$ws.Range("K122").Value = 9072
$ws.Range("M122").Value = -6622
# Row 126
$ws.Range("H126").Value = 2516.889
$ws.Range("I126").Value = 2206.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 6619.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -4149.5
$ws.Range("N126").Value = -19940
# Row 135
$ws.Range("H135").Value = 63790
$ws.Range("J135").Value = 63790
$ws.Range("L135").Value = 63790
$ws.Range("N135").Value = -73930
